# Change B5 value from 6 to 75
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 75

# Update the selection on the sheet to D12 (mirrors the cursor move saved in the file)
$ws.Range("D12").Select()
